$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "45-9=36"
$t.Cell(1, 2).Range.Text = "82+2=84"
$t.Cell(1, 3).Range.Text = "65-7=58"
$t.Cell(1, 4).Range.Text = "16+33=49"
$t.Cell(1, 5).Range.Text = "74-58=16"
$t.Cell(2, 1).Range.Text = "42-26=16"
$t.Cell(2, 2).Range.Text = "80-5=75"
$t.Cell(2, 3).Range.Text = "49+48=97"
$t.Cell(2, 4).Range.Text = "19-15=4"
$t.Cell(2, 5).Range.Text = "91-63=28"
$t.Cell(3, 1).Range.Text = "15+79=94"
$t.Cell(3, 2).Range.Text = "80-32=48"
$t.Cell(3, 3).Range.Text = "27+64=91"
$t.Cell(3, 4).Range.Text = "50+35=85"
$t.Cell(3, 5).Range.Text = "2+22=24"
$t.Cell(4, 1).Range.Text = "66+28=94"
$t.Cell(4, 2).Range.Text = "58-16=42"
$t.Cell(4, 3).Range.Text = "91-28=63"
$t.Cell(4, 4).Range.Text = "59+3=62"
$t.Cell(4, 5).Range.Text = "33-20=13"
$t.Cell(5, 1).Range.Text = "23+54=77"
$t.Cell(5, 2).Range.Text = "40+5=45"
$t.Cell(5, 3).Range.Text = "79-10=69"
$t.Cell(5, 4).Range.Text = "96-68=28"
$t.Cell(5, 5).Range.Text = "88+2=90"
$t.Cell(6, 1).Range.Text = "17+66=83"
$t.Cell(6, 2).Range.Text = "60-23=37"
$t.Cell(6, 3).Range.Text = "32+21=53"
$t.Cell(6, 4).Range.Text = "74-2=72"
$t.Cell(6, 5).Range.Text = "44+41=85"
$t.Cell(7, 1).Range.Text = "24+10=34"
$t.Cell(7, 2).Range.Text = "70-19=51"
$t.Cell(7, 3).Range.Text = "1+74=75"
$t.Cell(7, 4).Range.Text = "75+11=86"
$t.Cell(7, 5).Range.Text = "98-20=78"
$t.Cell(8, 1).Range.Text = "96-17=79"
$t.Cell(8, 2).Range.Text = "63-25=38"
$t.Cell(8, 3).Range.Text = "16-2=14"
$t.Cell(8, 4).Range.Text = "90-47=43"
$t.Cell(8, 5).Range.Text = "47+5=52"
$t.Cell(9, 1).Range.Text = "95-82=13"
$t.Cell(9, 2).Range.Text = "1+86=87"
$t.Cell(9, 3).Range.Text = "53-4=49"
$t.Cell(9, 4).Range.Text = "2+10=12"
$t.Cell(9, 5).Range.Text = "13+19=32"
$t.Cell(10, 1).Range.Text = "14+41=55"
$t.Cell(10, 2).Range.Text = "85-83=2"
$t.Cell(10, 3).Range.Text = "84+4=88"
$t.Cell(10, 4).Range.Text = "35+44=79"
$t.Cell(10, 5).Range.Text = "35+47=82"
$t.Cell(11, 1).Range.Text = "56-42=14"
$t.Cell(11, 2).Range.Text = "33+36=69"
$t.Cell(11, 3).Range.Text = "84-68=16"
$t.Cell(11, 4).Range.Text = "24+29=53"
$t.Cell(11, 5).Range.Text = "21+58=79"
$t.Cell(12, 1).Range.Text = "15-13=2"
$t.Cell(12, 2).Range.Text = "2+20=22"
$t.Cell(12, 3).Range.Text = "31+6=37"
$t.Cell(12, 4).Range.Text = "21+42=63"
$t.Cell(12, 5).Range.Text = "12+49=61"
$t.Cell(13, 1).Range.Text = "46+3=49"
$t.Cell(13, 2).Range.Text = "59-23=36"
$t.Cell(13, 3).Range.Text = "80-18=62"
$t.Cell(13, 4).Range.Text = "29-22=7"
$t.Cell(13, 5).Range.Text = "89-26=63"
$t.Cell(14, 1).Range.Text = "27+14=41"
$t.Cell(14, 2).Range.Text = "34-7=27"
$t.Cell(14, 3).Range.Text = "24+32=56"
$t.Cell(14, 4).Range.Text = "67+6=73"
$t.Cell(14, 5).Range.Text = "79+8=87"
$t.Cell(15, 1).Range.Text = "91-18=73"
$t.Cell(15, 2).Range.Text = "59-21=38"
$t.Cell(15, 3).Range.Text = "91-68=23"
$t.Cell(15, 4).Range.Text = "39-17=22"
$t.Cell(15, 5).Range.Text = "78-16=62"
$t.Cell(16, 1).Range.Text = "33+41=74"
$t.Cell(16, 2).Range.Text = "32+36=68"
$t.Cell(16, 3).Range.Text = "86-20=66"
$t.Cell(16, 4).Range.Text = "81-14=67"
$t.Cell(16, 5).Range.Text = "93-42=51"
$t.Cell(17, 1).Range.Text = "67-28=39"
$t.Cell(17, 2).Range.Text = "12+49=61"
$t.Cell(17, 3).Range.Text = "41+33=74"
$t.Cell(17, 4).Range.Text = "30+52=82"
$t.Cell(17, 5).Range.Text = "5+47=52"
$t.Cell(18, 1).Range.Text = "68+25=93"
$t.Cell(18, 2).Range.Text = "73+7=80"
$t.Cell(18, 3).Range.Text = "51+30=81"
$t.Cell(18, 4).Range.Text = "13+72=85"
$t.Cell(18, 5).Range.Text = "50-46=4"
$t.Cell(19, 1).Range.Text = "33+21=54"
$t.Cell(19, 2).Range.Text = "77-69=8"
$t.Cell(19, 3).Range.Text = "58+34=92"
$t.Cell(19, 4).Range.Text = "49-44=5"
$t.Cell(19, 5).Range.Text = "14+55=69"
$t.Cell(20, 1).Range.Text = "64+7=71"
$t.Cell(20, 2).Range.Text = "42+56=98"
$t.Cell(20, 3).Range.Text = "15+61=76"
$t.Cell(20, 4).Range.Text = "83+10=93"
$t.Cell(20, 5).Range.Text = "63-9=54"
